$p = $ppt.ActivePresentation

# ---- Slide 5: "Liveness and Performance Monitoring of SR Policy" ----
$s5 = $p.Slides.Item(5)

# Title: "Liveness and Performance Monitoring of SR Policy"
#     -> "Performance and Liveness Monitoring of SR Policy"
$s5Title = $s5.Shapes.Item(1)
$s5Title.TextFrame.TextRange.Text = "Performance and Liveness Monitoring of SR Policy"

# Content placeholder, first bullet:
# "Use PM probes (TWAMP Light/STAMP delay measurement messages) in Loopback Mode"
#     -> "Using PM probes (TWAMP Light/STAMP delay measurement messages) in Loopback Mode"
$s5Body = $s5.Shapes.Item(3)
$s5Body.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Using PM probes (TWAMP Light/STAMP delay measurement messages) in Loopback Mode"

# ---- Slide 6: "Enhanced Liveness and Performance Monitoring of SR Policy" ----
$s6 = $p.Slides.Item(6)

# Title: "Enhanced Liveness and Performance Monitoring of SR Policy"
#     -> "Enhanced Performance and Liveness Monitoring of SR Policy"
$s6Title = $s6.Shapes.Item(1)
$s6Title.TextFrame.TextRange.Text = "Enhanced Performance and Liveness Monitoring of SR Policy"

$s6Body = $s6.Shapes.Item(3)

# Shrink the content placeholder height from 2286000 EMU (180pt) to 2138362 EMU (168.375pt)
$s6Body.Height = 168.375

# First bullet: "Use PM probes in loopback mode enabled with network programming function"
#     -> "Using PM probes in loopback mode enabled with network programming function"
$s6Body.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Using PM probes in loopback mode enabled with network programming function"

# Third bullet: "The reflector node adds the receive timestamp in the payload of the received TWAMP Light or STAMP probe message without punting the probe message"
#     -> "Reflector node adds the receive timestamp in the payload of the received probe message without punting the message"
$s6Body.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Reflector node adds the receive timestamp in the payload of the received probe message without punting the message"
